$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove 3 rows from the data table (keeps the specially-bordered last
# table row intact, just shifting it up together with the footer rows).
$ws.Range("20:22").Delete()

# --- Header figures -------------------------------------------------
$ws.Range("E11").Value = 242153
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 7

# --- Data table (rows 16-23) ----------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "92255665"
$ws.Range("D16").Value = "OCTAVIO AUGUSTO MONTES HOYOS"
$ws.Range("E16").Value = "1704"
$ws.Range("F16").Value = 27578
$ws.Range("G16").Value = 800000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "92255665"
$ws.Range("D17").Value = "OCTAVIO AUGUSTO MONTES HOYOS"
$ws.Range("E17").Value = "1704"
$ws.Range("F17").Value = 29509
$ws.Range("G17").Value = 781242

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "8641226"
$ws.Range("D18").Value = "FABIAN ALBERTO SALAS CORONADO"
$ws.Range("E18").Value = "1908"
$ws.Range("F18").Value = 9600
$ws.Range("G18").Value = 1200000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "104300053"
$ws.Range("D19").Value = "CARLOS ENRIQUE PELAEZ AVILA"
$ws.Range("E19").Value = "1908"
$ws.Range("F19").Value = 5333
$ws.Range("G19").Value = 1000000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "104300053"
$ws.Range("D20").Value = "CARLOS ENRIQUE PELAEZ AVILA"
$ws.Range("E20").Value = "1909"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 1000000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "104300053"
$ws.Range("D21").Value = "CARLOS ENRIQUE PELAEZ AVILA"
$ws.Range("E21").Value = "1910"
$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 1000000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "9097544"
$ws.Range("D22").Value = "JOHN HAROL SEPULVEDA ZABALETA"
$ws.Range("E22").Value = "2102"
$ws.Range("F22").Value = 38133
$ws.Range("G22").Value = 1300000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "9097544"
$ws.Range("D23").Value = "JOHN HAROL SEPULVEDA ZABALETA"
$ws.Range("E23").Value = "2103"
$ws.Range("F23").Value = 52000
$ws.Range("G23").Value = 1300000
